$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 625, shifting the existing rows
# 625..661 down to 627..663.
$ws.Rows.Item(625).EntireRow.Insert() | Out-Null
$ws.Rows.Item(625).EntireRow.Insert() | Out-Null

# Populate the first new row (625) with the new weekly record.
$ws.Cells.Item(625, 1).Value = 11
$ws.Cells.Item(625, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(625, 3).Value = "Bíobío"
$ws.Cells.Item(625, 4).Value = 44931
$ws.Cells.Item(625, 5).Value = 8
$ws.Cells.Item(625, 6).Value = "Fruta"
$ws.Cells.Item(625, 7).Value = 100102
$ws.Cells.Item(625, 8).Value = "Cítricos"
$ws.Cells.Item(625, 9).Value = 100102003
$ws.Cells.Item(625, 10).Value = "Limón"
$ws.Cells.Item(625, 11).Value = "Sin especificar"
$ws.Cells.Item(625, 12).Value = "1a amarillo"
$ws.Cells.Item(625, 13).Value = 600
$ws.Cells.Item(625, 14).Value = 15000
$ws.Cells.Item(625, 15).Value = 16000
$ws.Cells.Item(625, 16).Value = 15500
$ws.Cells.Item(625, 17).Value = "$/malla 18 kilos"
$ws.Cells.Item(625, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(625, 19).Value = 861
$ws.Cells.Item(625, 20).Value = 18

# Populate the second new row (626) with the new weekly record.
$ws.Cells.Item(626, 1).Value = 11
$ws.Cells.Item(626, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(626, 3).Value = "Bíobío"
$ws.Cells.Item(626, 4).Value = 44931
$ws.Cells.Item(626, 5).Value = 8
$ws.Cells.Item(626, 6).Value = "Fruta"
$ws.Cells.Item(626, 7).Value = 100102
$ws.Cells.Item(626, 8).Value = "Cítricos"
$ws.Cells.Item(626, 9).Value = 100102003
$ws.Cells.Item(626, 10).Value = "Limón"
$ws.Cells.Item(626, 11).Value = "Sin especificar"
$ws.Cells.Item(626, 12).Value = "2a amarillo"
$ws.Cells.Item(626, 13).Value = 300
$ws.Cells.Item(626, 14).Value = 13000
$ws.Cells.Item(626, 15).Value = 13000
$ws.Cells.Item(626, 16).Value = 13000
$ws.Cells.Item(626, 17).Value = "$/malla 18 kilos"
$ws.Cells.Item(626, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(626, 19).Value = 722
$ws.Cells.Item(626, 20).Value = 18
